$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$ws1 = $wb.Worksheets.Item(1)

# Version: 1.0.0 -> 1.0.1
$ws1.Range("B3").Value = "1.0.1"

# Contact: "No display for ContactDetail" -> "MedCom (http://www.medcom.dk)"
$ws1.Range("B10").Value = "MedCom (http://www.medcom.dk)"

# Make room for a new "Jurisdiction" row right after "Contact" (row 10) by
# shifting the existing rows 11-14 (Description/Purpose/Copyright/Immutable)
# down to 12-15. Done via direct value moves (rather than Rows.Insert) so the
# shifted-in cells keep the same cell style as before instead of picking up a
# brand new default style.
$ws1.Range("A14:B14").Copy()
$ws1.Range("A15:B15").PasteSpecial(-4122)
$ws1.Range("A15").Value = $ws1.Range("A14").Value()
$ws1.Range("B15").Value = $ws1.Range("B14").Value()

$ws1.Range("A14").Value = $ws1.Range("A13").Value()
$ws1.Range("B14").Value = $ws1.Range("B13").Value()

$ws1.Range("A13").Value = $ws1.Range("A12").Value()
$ws1.Range("B13").Value = $ws1.Range("B12").Value()

$ws1.Range("A12").Value = $ws1.Range("A11").Value()
$ws1.Range("B12").Value = $ws1.Range("B11").Value()

$ws1.Range("A11").Value = "Jurisdiction"
$ws1.Range("B11").Value = ""

# --- Rename the second sheet ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Include #0"
